$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: after the run containing "address" (part of the
# "{d.applicant_address}" merge-field placeholder), insert a new run with
# text ":convCRLF" using the same BC Sans / 18-half-point formatting.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("{d.applicant_address}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not locate '{d.applicant_address}' placeholder"
}

$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="606C458E" w14:textId="4DFD77A8" w:rsidR="005035E1" w:rsidRPr="00934F2C" w:rsidRDefault="005035E1" w:rsidP="005035E1"><w:r w:rsidRPr="00934F2C"><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00934F2C"><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>d.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>applicant_</w:t></w:r><w:r w:rsidR="002059DB"><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>address</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>:convCRLF</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00934F2C"><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Change 2: remove the stray <w:lastRenderedPageBreak/> that precedes the
# "${d.cfrfee.feedata.estimatedtotaldue:formatN(2)}" merge field.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("estimatedtotaldue", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not locate 'estimatedtotaldue' merge field"
}
$para2 = $rng2.Paragraphs(1).Range

$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="49BB17BF" w14:textId="0DA71BAB" w:rsidR="005B49C3" w:rsidRDefault="005B49C3" w:rsidP="00FE533F"><w:pPr><w:tabs><w:tab w:val="left" w:pos="5575"/></w:tabs><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:b/><w:color w:val="0A3266"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="003327DF"><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:b/><w:color w:val="0A3266"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="003327DF"><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:b/><w:color w:val="0A3266"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>d.</w:t></w:r><w:r w:rsidR="0089382A" w:rsidRPr="003327DF"><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:b/><w:color w:val="0A3266"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>cfrfee</w:t></w:r><w:r w:rsidRPr="003327DF"><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:b/><w:color w:val="0A3266"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>.feedata.</w:t></w:r><w:r w:rsidR="00227924"><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:b/><w:color w:val="0A3266"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>estimatedtotaldue</w:t></w:r><w:r w:rsidR="004128DD" w:rsidRPr="00E75DC3"><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:b/><w:color w:val="0A3266"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>:formatN</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004128DD" w:rsidRPr="00E75DC3"><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:b/><w:color w:val="0A3266"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>(2)</w:t></w:r><w:r w:rsidRPr="003327DF"><w:rPr><w:rFonts w:ascii="BC Sans" w:eastAsia="BC Sans" w:hAnsi="BC Sans" w:cs="BC Sans"/><w:b/><w:color w:val="0A3266"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para2.InsertXML($xml2)

Write-Output "Both edits applied successfully"
